# "cambio de fracciones e historico" - update the quarterly report (a69_f26)
# to the next reporting period and correct the responsible-area name.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")

# --- Responsible area name correction (column AA, row 8) ---
$ws.Cells.Item(8, 27).Value = "Subdirección de Recursos Financieros (UPP)"

# --- Reporting period dates (row 8) ---
# Fecha de inicio del periodo que se informa
$ws.Cells.Item(8, 2).Value = 44743
# Fecha de término del periodo que se informa
$ws.Cells.Item(8, 3).Value = 44834

# Fecha de firma entrega de recursos (T8) gets date formatting applied
# (matches the other date columns, e.g. B8) even though left blank.
$ws.Cells.Item(8, 20).NumberFormat = $ws.Cells.Item(8, 2).NumberFormat

# Fecha de validación / Fecha de actualización (AB8 / AC8)
$ws.Cells.Item(8, 28).Value = 44844
$ws.Cells.Item(8, 29).Value = 44844

# --- Row heights ---
$ws.Rows.Item(3).RowHeight = 76.5
$ws.Rows.Item(8).RowHeight = 70.5

# --- Column AD width (Nota) ---
$ws.Columns.Item(30).ColumnWidth = 75.88

# --- Selection moved to A13 ---
$ws.Range("A13").Select() | Out-Null
